$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Kel'el Ware" (row 12) is removed from the roster; the players that followed
# shift up by one row, "Kessler Edwards" is added as a new entry, and
# "Michael Porter Jr." ends up as the final row.
$ws.Range("A12").Value = "Kristaps Porzingis"
$ws.Range("B12").Value = "PF,C"
$ws.Range("C12").Value = "Boston Celtics"

$ws.Range("A13").Value = "Naji Marshall"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Dallas Mavericks"

$ws.Range("A14").Value = "Alperen Sengün"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = "Houston Rockets"

$ws.Range("A15").Value = "Kessler Edwards"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Dallas Mavericks"

$ws.Range("A16").Value = "Julius Randle"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Minnesota Timberwolves"

$ws.Range("A17").Value = "Cam Thomas"
$ws.Range("B17").Value = "SG,SF"
$ws.Range("C17").Value = "Brooklyn Nets"

$ws.Range("A18").Value = "Michael Porter Jr."
$ws.Range("B18").Value = "SF,PF"
$ws.Range("C18").Value = "Denver Nuggets"
